# Update cryptocurrency price and volume(1h) values in the worksheet
# to reflect the latest data pulled on Mon May 15 18:50:28 UTC 2023.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.650.06"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.844.07"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("E4").Value = "  -1.72%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.14"
$ws.Range("E5").Value = "  -1.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.007"
$ws.Range("E6").Value = "  -1.73%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4295"
$ws.Range("E7").Value = "  -1.71%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3740"
$ws.Range("E8").Value = "  -1.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07315"
$ws.Range("E9").Value = "  -0.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8725"
$ws.Range("E10").Value = "  -0.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.45"
$ws.Range("E11").Value = "  -0.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.845.74"
$ws.Range("E12").Value = "  -0.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.713"
$ws.Range("E13").Value = "  +0.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.411"
$ws.Range("E14").Value = "  -1.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07095"
$ws.Range("E15").Value = "  -0.38%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.66"
$ws.Range("E16").Value = "  +4.38%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.011"
$ws.Range("E17").Value = "  -2.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008992"
$ws.Range("E18").Value = "  -0.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.008"
$ws.Range("E19").Value = "  -1.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.34"
$ws.Range("E20").Value = "  -0.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.651.87"
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.196"
$ws.Range("E22").Value = "  -1.65%  "
$ws.Range("E23").Value = "  -2.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.072.03"
$ws.Range("E24").Value = "  -0.77%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.969"
$ws.Range("E25").Value = "  -3.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.71"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.51"
$ws.Range("E27").Value = "  -0.78%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.156"
$ws.Range("E28").Value = "  +7.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.331"
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.77"
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08908"
$ws.Range("E31").Value = "  -1.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.220"
$ws.Range("E32").Value = "  +1.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7741"
$ws.Range("E33").Value = "  +0.69%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.525"
$ws.Range("E34").Value = "  -0.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.895"
$ws.Range("E35").Value = "  -2.84%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.009"
$ws.Range("E36").Value = "  -1.63%  "
$ws.Range("E37").Value = "  -1.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01975"
$ws.Range("E38").Value = "  +0.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05308"
$ws.Range("E39").Value = "  +0.81%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.882"
$ws.Range("E40").Value = "  +1.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.125"
$ws.Range("E41").Value = "  +4.43%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1688"
$ws.Range("E42").Value = "  +1.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5117"
$ws.Range("E43").Value = "  -0.99%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.766"
$ws.Range("E44").Value = "  +0.15%  "
$ws.Range("E45").Value = "  +0.64%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "107.21"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4741"
$ws.Range("E47").Value = "  +1.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06447"
$ws.Range("E48").Value = "  -2.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.008"
$ws.Range("E49").Value = "  -1.75%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.684"
$ws.Range("E50").Value = "  -0.83%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.840"
$ws.Range("E51").Value = "  -2.48%  "
